# Apply cryptocurrency price/volume updates for Sat Oct 14 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '26.901.19'
$cell.Style = 'Normal'
$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  -0.18%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.548.64'
$cell.Style = 'Normal'
$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  -0.42%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '206.49'
$cell.Style = 'Normal'
$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  -0.23%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  -0.39%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '22.18'
$cell.Style = 'Normal'
$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  +3.10%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  -0.86%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.0586'
$cell.Style = 'Normal'
$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  +0.53%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  -0.46%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '1.770.29'
$cell.Style = 'Normal'
$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  -0.44%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '1.549.23'
$cell.Style = 'Normal'
$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  -0.50%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  +0.71%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '0.517'
$cell.Style = 'Normal'
$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  +0.47%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '26.907.21'
$cell.Style = 'Normal'
$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  -0.16%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '61.55'
$cell.Style = 'Normal'
$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  -0.33%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '217.13'
$cell.Style = 'Normal'
$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  +1.22%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  +1.46%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  +0.03%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  -0.36%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  +0.03%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '9.23'
$cell.Style = 'Normal'
$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  +0.36%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  -0.59%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '154.20'
$cell.Style = 'Normal'
$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  +0.41%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '6.61'
$cell.Style = 'Normal'
$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  -0.76%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  +0.15%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  +0.51%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  -0.48%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '0.0466'
$cell.Style = 'Normal'
$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  +1.50%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '1.09'
$cell.Style = 'Normal'
$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  -0.83%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  -0.37%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '1.414.84'
$cell.Style = 'Normal'
$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  +3.14%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  +3.46%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  +1.96%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.964'
$cell.Style = 'Normal'
$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  -0.40%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  -0.11%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  +0.13%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  +0.69%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  -0.10%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  -0.39%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  +3.39%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  +3.60%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '0.998'
$cell.Style = 'Normal'
$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  +1.60%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '64.44'
$cell.Style = 'Normal'
$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  +1.25%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  +0.42%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '1.683.59'
$cell.Style = 'Normal'
$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  -0.46%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '87.42'
$cell.Style = 'Normal'
$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  +1.40%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  +4.16%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.0516'
$cell.Style = 'Normal'
$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  +1.84%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.0958'
$cell.Style = 'Normal'
$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  +0.25%  '
$cell.Style = 'Normal'
